$d = $word.ActiveDocument

# --- First new paragraph: hyperlink to the seaborn pairplot figure-size question ---
$r = $d.Content
$r.Collapse(0)
$r.InsertParagraphAfter()
$r.Collapse(0)
$h1 = $d.Hyperlinks.Add($r, "https://stackoverflow.com/questions/51400076/change-seaborn-pair-plot-figure-size")
$p1 = $d.Paragraphs.Last
$p1r = $p1.Range
$p1r.Characters.Item($p1r.Characters.Count).Text = " "

# --- Second new paragraph: hyperlink to the confusion matrix default figure-size question ---
$r2 = $d.Content
$r2.Collapse(0)
$r2.InsertParagraphAfter()
$r2.Collapse(0)
$h2 = $d.Hyperlinks.Add($r2, "https://stackoverflow.com/questions/61325314/how-to-change-plot-confusion-matrix-default-figure-size-in-sklearn-metrics-packa")
$p2 = $d.Paragraphs.Last
$p2r = $p2.Range
$p2r.Characters.Item($p2r.Characters.Count).Text = " "
